$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a paragraph whose (paragraph-mark-trimmed) text exactly
# matches the given string.
# ---------------------------------------------------------------------------
function Get-ParaByExactText($doc, $targetText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $targetText) {
            return $p
        }
    }
    return $null
}

# Helper: wrap a hand-built <w:p>...</w:p> (or run-only) fragment in the
# minimal WordOpenXML package envelope Range.InsertXML expects, then
# inject it into $rng (InsertXML *replaces* $rng's current contents).
function Set-RangeXml($rng, [string]$innerXml) {
    $pkg = '<?xml version="1.0"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    [void]$rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# Change 1: "Task 1." -> split "Task" | " 1." and drop the paragraph-mark
# run properties (<w:pPr><w:rPr><w:lang .../></w:rPr></w:pPr> disappears).
# Replace the *whole* paragraph (safe here - not the last paragraph in body)
# so the now-superfluous pPr goes away along with it.
# ---------------------------------------------------------------------------
$p1 = Get-ParaByExactText $d "Task 1."
if ($p1 -ne $null) {
    $xml1 = '<w:p>' +
              '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Task</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve"> 1.</w:t></w:r>' +
            '</w:p>'
    Set-RangeXml $p1.Range $xml1
}

# ---------------------------------------------------------------------------
# Change 2: "10. Who is Adam Hodge?" -> collapse the proof-errd runs
# ("10. Who " / "is" / " Adam " / "Hodge") into a single run, keep the
# trailing "?" run (lang=en-US), keep pPr untouched (re-emitted verbatim).
# ---------------------------------------------------------------------------
$p2 = Get-ParaByExactText $d "10. Who is Adam Hodge?"
if ($p2 -ne $null) {
    $xml2 = '<w:p>' +
              '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>10. Who is Adam Hodge</w:t></w:r>' +
              '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>?</w:t></w:r>' +
            '</w:p>'
    Set-RangeXml $p2.Range $xml2
}

# ---------------------------------------------------------------------------
# Change 3: "This text is too complicated ..." -> split into
# "This " | "situation" | " is too complicated to express my opinion,
# I would prefer to skip this question", all keeping lang=en-US.
# This is the LAST paragraph in the body - replacing Paragraph.Range
# (which includes the paragraph mark) would splice in an extra empty
# paragraph before sectPr, so only the visible-text sub-range is
# replaced here; pPr and the paragraph mark are left completely alone.
# ---------------------------------------------------------------------------
$target3 = "This text is too complicated to express my opinion, I would prefer to skip this question"
$p3 = Get-ParaByExactText $d $target3
if ($p3 -ne $null) {
    $textRange3 = $d.Range($p3.Range.Start, $p3.Range.Start + $target3.Length)
    $xml3 = '<w:p>' +
              '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">This </w:t></w:r>' +
              '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>situation</w:t></w:r>' +
              '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is too complicated to express my opinion, I would prefer to skip this question</w:t></w:r>' +
            '</w:p>'
    Set-RangeXml $textRange3 $xml3
}

Write-Output "OK"
